$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 3x3 cross-product data (ECs/FAPs/MuSCs sending x ECs/FAPs/MuSCs target) with recalculated TPM-based values
$data = @(
    @("ECs", "Ptn", "Ptprz1", "ECs", "1", "0.3333333333333333", "0.06520933333333333", "0.195628", "0.007116089623850998", "0.007116089623850999", "2", "0.6666666666666666", "0.06163433333333333", "0.184903", "0.006690894379667537", "0.006690894379667537", "0.00401913378711111", "0.03617220408399999", "4.761300406943512E-05", "4.761300406943512E-05"),
    @("ECs", "Ptn", "Ptprz1", "FAPs", "1", "0.3333333333333333", "0.06520933333333333", "0.195628", "0.007116089623850998", "0.007116089623850999", "3", "1", "0.1030763333333333", "0.309229", "0.01118975126488057", "0.01118975126488057", "0.00672153897911111", "0.060493850812", "7.962727286949023E-05", "7.962727286949024E-05"),
    @("ECs", "Ptn", "Ptprz1", "MuSCs", "1", "0.3333333333333333", "0.06520933333333333", "0.195628", "0.007116089623850998", "0.007116089623850999", "3", "1", "9.046962666666667", "27.140888", "0.9821193543554519", "0.9821193543554518", "0.5899464041848889", "5.309517637664", "0.006988849346912073", "0.006988849346912073"),
    @("FAPs", "Ptn", "Ptprz1", "ECs", "3", "1", "3.270036666666666", "9.81011", "0.3568488252184601", "0.3568488252184602", "2", "0.6666666666666666", "0.06163433333333333", "0.184903", "0.006690894379667537", "0.006690894379667537", "0.2015465299255555", "1.81391876933", "0.002387637799045158", "0.002387637799045158"),
    @("FAPs", "Ptn", "Ptprz1", "FAPs", "3", "1", "3.270036666666666", "9.81011", "0.3568488252184601", "0.3568488252184602", "3", "1", "0.1030763333333333", "0.309229", "0.01118975126488057", "0.01118975126488057", "0.3370633894655555", "3.03357050519", "0.00399304959335941", "0.003993049593359411"),
    @("FAPs", "Ptn", "Ptprz1", "MuSCs", "3", "1", "3.270036666666666", "9.81011", "0.3568488252184601", "0.3568488252184602", "3", "1", "9.046962666666667", "27.140888", "0.9821193543554519", "0.9821193543554518", "29.58389964196445", "266.25509677768", "0.3504681378260555", "0.3504681378260556"),
    @("MuSCs", "Ptn", "Ptprz1", "ECs", "3", "1", "5.828401", "17.485203", "0.6360350851576888", "0.6360350851576889", "2", "0.6666666666666666", "0.06163433333333333", "0.184903", "0.006690894379667537", "0.006690894379667537", "0.3592296100343333", "3.233066490309", "0.004255643576552943", "0.004255643576552944"),
    @("MuSCs", "Ptn", "Ptprz1", "FAPs", "3", "1", "5.828401", "17.485203", "0.6360350851576888", "0.6360350851576889", "3", "1", "0.1030763333333333", "0.309229", "0.01118975126488057", "0.01118975126488057", "0.6007702042763333", "5.406931838487", "0.007117074398651672", "0.007117074398651673"),
    @("MuSCs", "Ptn", "Ptprz1", "MuSCs", "3", "1", "5.828401", "17.485203", "0.6360350851576888", "0.6360350851576889", "3", "1", "9.046962666666667", "27.140888", "0.9821193543554519", "0.9821193543554518", "52.72932625336267", "474.563936280264", "0.6246623671824842", "0.6246623671824842"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        if ($c -lt 4) {
            $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
        } else {
            $ws.Cells.Item($r, $c + 1).Value2 = [double]$row[$c]
        }
    }
}